$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sampleIds = @(
    "C0024276",
    "C0024269",
    "C0024268",
    "C0024267",
    "C0024262",
    "C0024094",
    "C0024067",
    "C0024064",
    "C0023986",
    "C0023689",
    "C0023676",
    "C0023512",
    "C0023380",
    "C0023369",
    "C0023367",
    "C0023365",
    "C0023344",
    "C0023342",
    "C0023248",
    "C0023224",
    "C0023217",
    "C0021055",
    "C0020968",
    "C0020967",
    "C0020966",
    "C0020965",
    "C0020953",
    "C0020952"
)

$ws.Range("A61:K61").Copy()
$destRange = $ws.Range("A62:K89")
$destRange.PasteSpecial(-4122)

for ($i = 0; $i -lt $sampleIds.Count; $i++) {
    $row = 62 + $i
    $ws.Cells.Item($row, 1).Value2 = $sampleIds[$i]
    $ws.Cells.Item($row, 8).Value2 = "Ctrl"
}

$excel.CutCopyMode = 0

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:L89"))

# Select the first empty row below the newly added data (mirrors the
# author's final cursor position) and scroll the new rows into view.
$ws.Range("A90").Select()
try {
    $excel.ActiveWindow.ScrollRow = 56
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Window-level scroll position isn't always exposed; selection above
    # already captures the meaningful part of the view state.
}
